# Generate Report for Archive
#
# 1) The status text "Ready for handoff" becomes "In Translation" for the
#    localization entry (Overview!E2/F2, zh-cn!C2, de-de!C2 all share that
#    string).
# 2) The Status column narrows (its autofit width shrinks) on every sheet
#    that shows it: Overview columns E & F, and column C on the zh-cn /
#    de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Narrow the now-updated Status columns. 12.5 is the closest width this
# engine's column-width grid can represent to the target 13.41 chars.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
